# Auto-generated script applying updated market price data
# to the Halicarnassus_Profits workbook (Leve profit tracker).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 80
$ws.Range("H80").Value = 198.81818
$ws.Range("I80").Value = 85.28570999999999
$ws.Range("J80").Value = 397.5
$ws.Range("K80").Value = 255.85713
$ws.Range("L80").Value = 1192.5
$ws.Range("M80").Value = 742.14287
$ws.Range("N80").Value = -3188.5
# Row 83
$ws.Range("H83").Value = 198.81818
$ws.Range("I83").Value = 85.28570999999999
$ws.Range("J83").Value = 397.5
$ws.Range("K83").Value = 767.57139
$ws.Range("L83").Value = 3577.5
$ws.Range("M83").Value = 4224.42861
$ws.Range("N83").Value = -13561.5
# Row 86
$ws.Range("H86").Value = 1712.5
$ws.Range("I86").Value = 1100
$ws.Range("J86").Value = 1916.6666
$ws.Range("K86").Value = 1100
$ws.Range("L86").Value = 1916.6666
$ws.Range("M86").Value = 23
$ws.Range("N86").Value = -4162.6666
# Row 89
$ws.Range("H89").Value = 1712.5
$ws.Range("I89").Value = 1100
$ws.Range("J89").Value = 1916.6666
$ws.Range("K89").Value = 5500
$ws.Range("L89").Value = 9583.333000000001
$ws.Range("M89").Value = 116
$ws.Range("N89").Value = -20815.333
# Row 100
$ws.Range("H100").Value = 897.38464
$ws.Range("I100").Value = 883.4545000000001
$ws.Range("K100").Value = 883.4545000000001
$ws.Range("M100").Value = -342.4545000000001
# Row 113
$ws.Range("H113").Value = 5655
$ws.Range("I113").Value = 6017
$ws.Range("K113").Value = 6017
$ws.Range("M113").Value = -2763
# Row 130
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
# Row 137
$ws.Range("H137").Value = 1999.9375
$ws.Range("I137").Value = 916.2222
$ws.Range("J137").Value = 3393.2856
$ws.Range("K137").Value = 2748.6666
$ws.Range("L137").Value = 10179.8568
$ws.Range("M137").Value = -198.6666
$ws.Range("N137").Value = -15279.8568
# Row 138
$ws.Range("H138").Value = 2490.1428
$ws.Range("I138").Value = 1036.4166
$ws.Range("K138").Value = 3109.2498
$ws.Range("M138").Value = 2030.7502

$ws = $wb.Worksheets.Item("ARM")
# Row 19
$ws.Range("H19").Value = 15000
$ws.Range("I19").Value = 10000
$ws.Range("K19").Value = 10000
$ws.Range("M19").Value = -9771
# Row 97
$ws.Range("H97").Value = 1064.9166
$ws.Range("I97").Value = 1078.5
$ws.Range("K97").Value = 1078.5
$ws.Range("M97").Value = -582.5
# Row 110
$ws.Range("H110").Value = 2374.875
$ws.Range("I110").Value = 2247.6667
$ws.Range("K110").Value = 2247.6667
$ws.Range("M110").Value = -202.6667000000002
# Row 122
$ws.Range("H122").Value = 1996.6666
$ws.Range("I122").Value = 1997.5
$ws.Range("K122").Value = 5992.5
$ws.Range("M122").Value = -3542.5

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 5183.778
$ws.Range("I86").Value = 1288.5
$ws.Range("J86").Value = 8300
$ws.Range("K86").Value = 1288.5
$ws.Range("L86").Value = 8300
$ws.Range("M86").Value = -165.5
$ws.Range("N86").Value = -10546
# Row 89
$ws.Range("H89").Value = 5183.778
$ws.Range("I89").Value = 1288.5
$ws.Range("J89").Value = 8300
$ws.Range("K89").Value = 6442.5
$ws.Range("L89").Value = 41500
$ws.Range("M89").Value = -826.5
$ws.Range("N89").Value = -52732
# Row 94
$ws.Range("H94").Value = 909.75
$ws.Range("I94").Value = 909.75
$ws.Range("K94").Value = 909.75
$ws.Range("M94").Value = -458.75
# Row 105
$ws.Range("H105").Value = 1712.6875
$ws.Range("I105").Value = 1626.8667
$ws.Range("K105").Value = 1626.8667
$ws.Range("M105").Value = 120.1333
# Row 107
$ws.Range("H107").Value = 5008.3
$ws.Range("I107").Value = 1680.5
$ws.Range("K107").Value = 1680.5
$ws.Range("M107").Value = 239.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5491.2085
$ws.Range("I31").Value = 1470.4546
$ws.Range("K31").Value = 1470.4546
$ws.Range("M31").Value = -1175.4546
# Row 34
$ws.Range("H34").Value = 5491.2085
$ws.Range("I34").Value = 1470.4546
$ws.Range("K34").Value = 1470.4546
$ws.Range("M34").Value = -1268.4546
# Row 44
$ws.Range("H44").Value = 15900
$ws.Range("I44").Value = 15900
$ws.Range("K44").Value = 15900
$ws.Range("M44").Value = -15458
# Row 50
$ws.Range("H50").Value = 90357.336
$ws.Range("J50").Value = 90357.336
$ws.Range("L50").Value = 90357.336
$ws.Range("N50").Value = -91607.336
# Row 55
$ws.Range("H55").Value = 8000
$ws.Range("I55").Value = 8000
$ws.Range("K55").Value = 8000
$ws.Range("M55").Value = -7685
# Row 86
$ws.Range("H86").Value = 4250
$ws.Range("I86").Value = 3500
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 3500
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -2377
$ws.Range("N86").Value = -7246
# Row 89
$ws.Range("H89").Value = 4250
$ws.Range("I89").Value = 3500
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 17500
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -11884
$ws.Range("N89").Value = -36232
# Row 99
$ws.Range("H99").Value = 2435.3635
$ws.Range("I99").Value = 2489.4
$ws.Range("J99").Value = 1895
$ws.Range("K99").Value = 2489.4
$ws.Range("L99").Value = 1895
$ws.Range("M99").Value = -991.4000000000001
$ws.Range("N99").Value = -4891
# Row 107
$ws.Range("H107").Value = 486.56522
$ws.Range("I107").Value = 243.05882
$ws.Range("J107").Value = 1176.5
$ws.Range("K107").Value = 243.05882
$ws.Range("L107").Value = 1176.5
$ws.Range("M107").Value = 1676.94118
$ws.Range("N107").Value = -5016.5
# Row 112
$ws.Range("H112").Value = 45000
$ws.Range("J112").Value = 45000
$ws.Range("L112").Value = 45000
$ws.Range("N112").Value = -47954
# Row 126
$ws.Range("H126").Value = 2435.3635
$ws.Range("I126").Value = 2489.4
$ws.Range("J126").Value = 1895
$ws.Range("K126").Value = 7468.200000000001
$ws.Range("L126").Value = 5685
$ws.Range("M126").Value = -4998.200000000001
$ws.Range("N126").Value = -10625

$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 117.85714
$ws.Range("I23").Value = 61.25
$ws.Range("J23").Value = 193.33333
$ws.Range("K23").Value = 183.75
$ws.Range("L23").Value = 579.99999
$ws.Range("M23").Value = 51.25
$ws.Range("N23").Value = -1049.99999
# Row 38
$ws.Range("H38").Value = 467.2381
$ws.Range("I38").Value = 441.82352
$ws.Range("K38").Value = 1325.47056
$ws.Range("M38").Value = -978.47056
# Row 98
$ws.Range("H98").Value = 183.33333
$ws.Range("I98").Value = 190
$ws.Range("J98").Value = 170
$ws.Range("K98").Value = 570
$ws.Range("L98").Value = 510
$ws.Range("M98").Value = 928
$ws.Range("N98").Value = -3506
# Row 113
$ws.Range("H113").Value = 1349.2858
$ws.Range("I113").Value = 930
$ws.Range("K113").Value = 2790
$ws.Range("M113").Value = -620
# Row 122
$ws.Range("H122").Value = 750
$ws.Range("J122").Value = 750
$ws.Range("L122").Value = 6750
$ws.Range("N122").Value = -11650
# Row 132
$ws.Range("H132").Value = 4217
$ws.Range("I132").Value = 3500
$ws.Range("J132").Value = 4421.857
$ws.Range("K132").Value = 31500
$ws.Range("L132").Value = 39796.713
$ws.Range("M132").Value = -28970
$ws.Range("N132").Value = -44856.713

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1801.68
$ws.Range("I102").Value = 1361.4546
$ws.Range("K102").Value = 1361.4546
$ws.Range("M102").Value = 260.5454
# Row 122
$ws.Range("H122").Value = 3138.6667
$ws.Range("I122").Value = 3016.1667
$ws.Range("K122").Value = 9048.500100000001
$ws.Range("M122").Value = -6598.500100000001
# Row 126
$ws.Range("H126").Value = 4959.6665
$ws.Range("I126").Value = 2489.3333
$ws.Range("K126").Value = 7467.999899999999
$ws.Range("M126").Value = -4997.999899999999
# Row 132
$ws.Range("H132").Value = 3174.6
$ws.Range("I132").Value = 2843.5
$ws.Range("J132").Value = 4499
$ws.Range("K132").Value = 8530.5
$ws.Range("L132").Value = 13497
$ws.Range("M132").Value = -6000.5
$ws.Range("N132").Value = -18557

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 2910.4
$ws.Range("I122").Value = 2844.889
$ws.Range("K122").Value = 8534.667000000001
$ws.Range("M122").Value = -6084.667000000001
# Row 132
$ws.Range("H132").Value = 4888
$ws.Range("J132").Value = 4875
$ws.Range("L132").Value = 14625
$ws.Range("N132").Value = -19685

$ws = $wb.Worksheets.Item("WVR")
# Row 112
$ws.Range("H112").Value = 17856
$ws.Range("J112").Value = 17856
$ws.Range("L112").Value = 17856
$ws.Range("N112").Value = -20810
# Row 122
$ws.Range("H122").Value = 1316.8823
$ws.Range("I122").Value = 1180.7273
$ws.Range("J122").Value = 1566.5
$ws.Range("K122").Value = 3542.1819
$ws.Range("L122").Value = 4699.5
$ws.Range("M122").Value = -1092.1819
$ws.Range("N122").Value = -9599.5
# Row 132
$ws.Range("H132").Value = 1951.6522
$ws.Range("I132").Value = 1858.5454
$ws.Range("K132").Value = 5575.6362
$ws.Range("M132").Value = -3045.6362

